$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns hold numeric-looking text
# (e.g. "16.09", "4.60", "39.469.99") that must stay exact text rather
# than being re-interpreted as numbers. Temporarily mark the cell as
# Text before assigning, then restore the default "Normal" style so the
# cell keeps no explicit style index (matching the original workbook).
$textCells = @(
    'D2',
    'E2',
    'D3',
    'E3',
    'E4',
    'D5',
    'E5',
    'E6',
    'D7',
    'E7',
    'E8',
    'E9',
    'E10',
    'E11',
    'D12',
    'E12',
    'D13',
    'E13',
    'D14',
    'E14',
    'E15',
    'D16',
    'E16',
    'D17',
    'E17',
    'D18',
    'E18',
    'D19',
    'E19',
    'D20',
    'E20',
    'D21',
    'E21',
    'D22',
    'E22',
    'E23',
    'D24',
    'E24',
    'E25',
    'D26',
    'E26',
    'D27',
    'E27',
    'E28',
    'E29',
    'D30',
    'E30',
    'D31',
    'E31',
    'E32',
    'D33',
    'E33',
    'D34',
    'E34',
    'D35',
    'E35',
    'E36',
    'D37',
    'E37',
    'D38',
    'E38',
    'E39',
    'D40',
    'E40',
    'E41',
    'D42',
    'E42',
    'D43',
    'E43',
    'E44',
    'D45',
    'E45',
    'D46',
    'E46',
    'E47',
    'E48',
    'D49',
    'E49',
    'D50',
    'E50',
    'D51',
    'E51',
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '39.469.99'
$ws.Range('E2').Value = '  +1.87%  '
$ws.Range('D3').Value = '2.160.23'
$ws.Range('E3').Value = '  +2.96%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '228.11'
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('D7').Value = '64.08'
$ws.Range('E7').Value = '  +2.95%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +2.24%  '
$ws.Range('E10').Value = '  +1.68%  '
$ws.Range('E11').Value = '  +0.19%  '
$ws.Range('D12').Value = '16.09'
$ws.Range('E12').Value = '  +1.62%  '
$ws.Range('D13').Value = '2.479.34'
$ws.Range('E13').Value = '  +2.85%  '
$ws.Range('D14').Value = '22.14'
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').Value = '5.53'
$ws.Range('E16').Value = '  +0.55%  '
$ws.Range('D17').Value = '2.140.83'
$ws.Range('E17').Value = '  +1.99%  '
$ws.Range('D18').Value = '39.399.22'
$ws.Range('E18').Value = '  +1.64%  '
$ws.Range('D19').Value = '71.87'
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('D20').Value = '6.11'
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('D21').Value = '0.0₃0852'
$ws.Range('E21').Value = '  +1.61%  '
$ws.Range('D22').Value = '229.93'
$ws.Range('E22').Value = '  +0.90%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').Value = '2.36'
$ws.Range('E24').Value = '  +1.17%  '
$ws.Range('E25').Value = '  -2.76%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '172.60'
$ws.Range('E26').Value = '  +0.35%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '9.52'
$ws.Range('E27').Value = '  -0.63%  '
$ws.Range('E28').Value = '  +1.78%  '
$ws.Range('E29').Value = '  +2.68%  '
$ws.Range('D30').Value = '1.42'
$ws.Range('E30').Value = '  +0.29%  '
$ws.Range('D31').Value = '2.58'
$ws.Range('E31').Value = '  +4.60%  '
$ws.Range('E32').Value = '  +1.00%  '
$ws.Range('D33').Value = '4.60'
$ws.Range('E33').Value = '  +1.54%  '
$ws.Range('D34').Value = '7.14'
$ws.Range('E34').Value = '  +8.50%  '
$ws.Range('D35').Value = '4.74'
$ws.Range('E35').Value = '  -0.25%  '
$ws.Range('E36').Value = '  -0.64%  '
$ws.Range('D37').Value = '2.43'
$ws.Range('E37').Value = '  +0.78%  '
$ws.Range('D38').Value = '3.56'
$ws.Range('E38').Value = '  +0.27%  '
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').Value = '103.43'
$ws.Range('E40').Value = '  +1.23%  '
$ws.Range('E41').Value = '  +0.77%  '
$ws.Range('D42').Value = '17.69'
$ws.Range('E42').Value = '  -3.41%  '
$ws.Range('D43').Value = '1.531.99'
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('E44').Value = '  +3.78%  '
$ws.Range('D45').Value = '4.32'
$ws.Range('E45').Value = '  +4.43%  '
$ws.Range('D46').Value = '0.0932'
$ws.Range('E46').Value = '  +2.27%  '
$ws.Range('E47').Value = '  +0.55%  '
$ws.Range('E48').Value = '  +5.40%  '
$ws.Range('D49').Value = '7.82'
$ws.Range('E49').Value = '  +0.48%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.362.95'
$ws.Range('E50').Value = '  +2.85%  '
$ws.Range('B51').Value = 'Celestia'
$ws.Range('C51').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D51').Value = '9.09'
$ws.Range('E51').Value = '  +23.76%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
